# Actualización automática 2025-10-30 15:30:08
# Updates the sales figures for client "JARAMILLO CARVAJAL NICOLAS ESTEBAN"
# (advisor "HIDALGO HIDALGO PEDRO GUSTAVO") to reflect an additional
# 95.04 in sales for the "240X80 PORCELANATO" group during octubre, and
# propagates that change through the dependent summary/total cells.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": per-client sales broken down by product group
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D11").Value = 743.04

# --- Sheet "VENTA MENSUAL": per-client sales broken down by month
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F11").Value = 6951.07
$wsMensual.Range("F23").Value = 27491.2

# --- Sheet "CUMPLIMIENTO MENSUAL": budget compliance summary by group
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D3").Value = 1477.47
$wsCumplimiento.Range("E3").Value = 4027.14890386263
$wsCumplimiento.Range("F3").Value = 0.2684055019618613
$wsCumplimiento.Range("D14").Value = 27491.2
$wsCumplimiento.Range("E14").Value = 27933.5414788039
$wsCumplimiento.Range("F14").Value = 0.496009530518306
